$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A; this shifts the whole A:D block to B:E
# (also shifts the column width metadata the same way the real edit did).
$ws.Columns.Item(1).Insert()

# ---- Row 1 (headers): new A1 = "ID"; the rest stayed shifted from the insert ----
$ws.Range("A1").Value = "ID"

# ---- Data rows: put the DataID number back in column A (undo the shift there) ----
# and put the new building_* identifier into column B (shifting Name/Pollution/Cost
# from B:D back out to C:E, matching the original, non-inserted layout for those rows).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "building_towncenter"
$ws.Range("C2").Value = "Town Center"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "building_farm"
$ws.Range("C3").Value = "Farm"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "building_factory"
$ws.Range("C4").Value = "Factory"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "building_filterationplant"
$ws.Range("C5").Value = "Filteration Plant"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "building_house"
$ws.Range("C6").Value = "House"

# Match the final cell selection left in the saved workbook.
$null = $ws.Range("D4").Select()
